$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: B1 "Task" -> "Assign Date" ---
$ws.Range("B1").Value = "Assign Date"

# --- Row 2: new cells B2:F2 ---
# Write in an order that reproduces the shared-string table order of the
# target file (Done, no delay, then the long Django description).
$ws.Range("D2").Value = "Done"
$ws.Range("F2").Value = "no delay"
$ws.Range("E2").Value = "How to install the Django and its command lines:- first you should install the python then you should install the django environment first make directory `$mkdir django then goto the django> python -m venv myvenv the activate command django>myenv\Scripts\activate  then django>python -m pip install Django  then make your first project django>django-admin.exe startproject project_name then do changes in setting.py files like time installed apps etc. then make the app inside your project like  django>project_name>python manage.py startapp app_name  then command for run the server is >python manage.py runserver"

# Dates as serials so no stray time-of-day fraction is attached.
$ws.Range("B2").Value = 44161
$ws.Range("C2").Value = 44162
$ws.Range("B2:C2").NumberFormat = "mm-dd-yy"

# Alignment for the description cell and row-13 marker cell.
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4160
$ws.Range("E13").HorizontalAlignment = -4108

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 40.05
$ws.Rows.Item(13).RowHeight = 15.75

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 11.5546875
$ws.Columns.Item(2).ColumnWidth = 10.33203125
$ws.Columns.Item(3).ColumnWidth = 10.33203125
$ws.Columns.Item(4).ColumnWidth = 5.33203125
$ws.Columns.Item(5).ColumnWidth = 67.33203125
$ws.Columns.Item(6).ColumnWidth = 7.88671875

# --- Selection matches the saved view in the target file ---
$ws.Range("E10").Select()
